# Auto-generated Excel COM-interop edit script
# Applies scheduled market-data refresh values to the Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1327.2106
$ws.Range("I19").Value = 1291.2307
$ws.Range("K19").Value = 1291.2307
$ws.Range("M19").Value = -1116.2307

$ws.Range("H28").Value = 2413.6667
$ws.Range("I28").Value = 2413.6667
$ws.Range("K28").Value = 2413.6667
$ws.Range("M28").Value = -1928.6667

$ws.Range("H92").Value = 1335.25
$ws.Range("I92").Value = 1335.25
$ws.Range("K92").Value = 1335.25
$ws.Range("M92").Value = -87.25

$ws.Range("H132").Value = 1555
$ws.Range("I132").Value = 1555
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4665
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -2135
$ws.Range("M132").ClearContents()

$ws.Range("H137").Value = 8165.8887
$ws.Range("I137").Value = 59998
$ws.Range("K137").Value = 179994
$ws.Range("M137").Value = -177444

$ws.Range("H138").Value = 4015.2046
$ws.Range("I138").Value = 2904.75
$ws.Range("J138").Value = 4261.972
$ws.Range("K138").Value = 8714.25
$ws.Range("L138").Value = 12785.916
$ws.Range("M138").Value = -3574.25
$ws.Range("N138").Value = -23065.916

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3989.2058
$ws.Range("I32").Value = 3144.8125
$ws.Range("J32").Value = 17499.5
$ws.Range("K32").Value = 3144.8125
$ws.Range("L32").Value = 17499.5
$ws.Range("M32").Value = -2857.8125
$ws.Range("N32").Value = -18073.5

$ws.Range("H45").Value = 2471.8125
$ws.Range("I45").Value = 2200.4
$ws.Range("J45").Value = 2924.1667
$ws.Range("K45").Value = 2200.4
$ws.Range("L45").Value = 2924.1667
$ws.Range("M45").Value = -1823.4
$ws.Range("N45").Value = -3678.1667

$ws.Range("H61").Value = 1901.75
$ws.Range("I61").Value = 1791.4736
$ws.Range("J61").Value = 3997
$ws.Range("K61").Value = 1791.4736
$ws.Range("L61").Value = 3997
$ws.Range("M61").Value = -1579.4736
$ws.Range("N61").Value = -4421

$ws.Range("H74").Value = 500
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 500
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H132").Value = 3285.5625
$ws.Range("I132").Value = 2233.7273
$ws.Range("J132").Value = 5599.6
$ws.Range("K132").Value = 6701.1819
$ws.Range("L132").Value = 16798.8
$ws.Range("M132").Value = -4171.1819
$ws.Range("N132").Value = -21858.8

$ws.Range("H136").Value = 1901.75
$ws.Range("I136").Value = 1791.4736
$ws.Range("J136").Value = 3997
$ws.Range("K136").Value = 5374.4208
$ws.Range("L136").Value = 11991
$ws.Range("M136").Value = -2824.4208
$ws.Range("N136").Value = -17091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 802.0909
$ws.Range("I99").Value = 810.9
$ws.Range("J99").Value = 714
$ws.Range("K99").Value = 810.9
$ws.Range("L99").Value = 714
$ws.Range("M99").Value = 687.1
$ws.Range("N99").Value = -3710

$ws.Range("H141").Value = 125000
$ws.Range("J141").Value = 125000
$ws.Range("L141").Value = 125000
$ws.Range("N141").Value = -135360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6396.4287
$ws.Range("I86").Value = 6345
$ws.Range("J86").Value = 6489
$ws.Range("K86").Value = 6345
$ws.Range("L86").Value = 6489
$ws.Range("M86").Value = -5222
$ws.Range("N86").Value = -8735

$ws.Range("H89").Value = 6396.4287
$ws.Range("I89").Value = 6345
$ws.Range("J89").Value = 6489
$ws.Range("K89").Value = 31725
$ws.Range("L89").Value = 32445
$ws.Range("M89").Value = -26109
$ws.Range("N89").Value = -43677

$ws.Range("H134").Value = 7188.9
$ws.Range("I134").Value = 7543.3335
$ws.Range("K134").Value = 22630.0005
$ws.Range("M134").Value = -20095.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 189998
$ws.Range("J37").Value = 189998
$ws.Range("L37").Value = 569994
$ws.Range("N37").Value = -570218

$ws.Range("H80").Value = 8995.200000000001
$ws.Range("I80").Value = 8998.5
$ws.Range("J80").Value = 8993
$ws.Range("K80").Value = 26995.5
$ws.Range("L80").Value = 26979
$ws.Range("M80").Value = -26059.5
$ws.Range("N80").Value = -28851

$ws.Range("H83").Value = 8995.200000000001
$ws.Range("I83").Value = 8998.5
$ws.Range("J83").Value = 8993
$ws.Range("K83").Value = 80986.5
$ws.Range("L83").Value = 80937
$ws.Range("M83").Value = -76306.5
$ws.Range("N83").Value = -90297

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 298.5
$ws.Range("I107").Value = 298.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 298.5
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = 1621.5
$ws.Range("M107").ClearContents()

$ws.Range("H126").Value = 1954.9
$ws.Range("I126").Value = 1387.5
$ws.Range("K126").Value = 4162.5
$ws.Range("M126").Value = -1692.5

$ws.Range("H132").Value = 4326.6665
$ws.Range("I132").Value = 3992.2
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 11976.6
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -9446.599999999999
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3373.8333
$ws.Range("I7").Value = 2848.6
$ws.Range("K7").Value = 2848.6
$ws.Range("M7").Value = -2736.6

$ws.Range("H16").Value = 1424.6666
$ws.Range("I16").Value = 1424.6666
$ws.Range("K16").Value = 1424.6666
$ws.Range("M16").Value = -1254.6666

$ws.Range("H22").Value = 3900
$ws.Range("I22").Value = 3900
$ws.Range("K22").Value = 3900
$ws.Range("M22").Value = -3605

$ws.Range("H27").Value = 3900
$ws.Range("I27").Value = 3900
$ws.Range("K27").Value = 3900
$ws.Range("M27").Value = -3793

$ws.Range("H55").Value = 473.75
$ws.Range("I55").Value = 497.5
$ws.Range("K55").Value = 497.5
$ws.Range("M55").Value = -324.5

$ws.Range("H61").Value = 4830.75
$ws.Range("I61").Value = 4858
$ws.Range("J61").Value = 4749
$ws.Range("K61").Value = 4858
$ws.Range("L61").Value = 4749
$ws.Range("M61").Value = -4656
$ws.Range("N61").Value = -5153

$ws.Range("H93").Value = 562.6
$ws.Range("I93").Value = 453.5
$ws.Range("K93").Value = 453.5
$ws.Range("M93").Value = 794.5

$ws.Range("H113").Value = 4830.75
$ws.Range("I113").Value = 4858
$ws.Range("J113").Value = 4749
$ws.Range("K113").Value = 4858
$ws.Range("L113").Value = 4749
$ws.Range("M113").Value = -2688
$ws.Range("N113").Value = -9089

$ws.Range("H126").Value = 3373.8333
$ws.Range("I126").Value = 2848.6
$ws.Range("K126").Value = 8545.799999999999
$ws.Range("M126").Value = -6075.799999999999

$ws.Range("H132").Value = 4410.4116
$ws.Range("I132").Value = 3726.3635
$ws.Range("J132").Value = 5664.5
$ws.Range("K132").Value = 11179.0905
$ws.Range("L132").Value = 16993.5
$ws.Range("M132").Value = -8649.0905
$ws.Range("N132").Value = -22053.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 46500
$ws.Range("J115").Value = 46500
$ws.Range("L115").Value = 46500
$ws.Range("N115").Value = -49634

$ws.Range("H136").Value = 6750.154
$ws.Range("I136").Value = 6750.154
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 20250.462
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -17700.462
$ws.Range("M136").ClearContents()
